$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated trait values (kernel count image results output)
$data = @{
    2 = @{ D = 63.34496825573626; E = 0.9427239224483256; F = 5.299253731343284; G = 15.47462686567164; H = 159; I = 0.2724400683851932; J = 102; K = 0.2628502434054487 }
    3 = @{ D = 57.22118149922032; E = 0.9581351127950248; F = 4.956716417910448; G = 14.85;              H = 140; I = 0.2864291064666009; J = 96;  K = 0.2555510866841167 }
    4 = @{ D = 52.52100718422812; E = 0.90365981393659;   F = 4.452985074626866; G = 16.90522388059702; H = 149; I = 0.2070724859126186; J = 112; K = 0.1969933249530599 }
    5 = @{ D = 42.56708788148809; E = 0.91010710845327;   F = 4.452985074626866; G = 11.76716417910448; H = 106; I = 0.2920866546947045; J = 63;  K = 0.293265084492251 }
    6 = @{ D = 30.76062513922923; E = 0.7182461905987475; F = 3.707462686567164; G = 11.2634328358209;  H = 69;  I = 0.2442220960415315; J = 58;  K = 0.2420939835690309 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    $ws.Range("D$row").Value = $cols.D
    $ws.Range("E$row").Value = $cols.E
    $ws.Range("F$row").Value = $cols.F
    $ws.Range("G$row").Value = $cols.G
    $ws.Range("H$row").Value = $cols.H
    $ws.Range("I$row").Value = $cols.I
    $ws.Range("J$row").Value = $cols.J
    $ws.Range("K$row").Value = $cols.K
}
